# The commit swaps the presentation's active theme colours away from the
# "Integral" / "Red Violet" palette that ships in ppt/theme/theme1.xml
# toward the standard Office default palette (the one originally parked,
# unused, in ppt/theme/theme2.xml as "Office Theme" / "Office").
#
# theme1.xml is the theme actually wired to the slide master (and to the
# presentation part itself), so recolouring it via the presentation's
# ThemeColorScheme is what makes the deck visually become the "Office"
# design. PowerPoint's 12 theme colour slots map, in order, to:
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
#   8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# RGB() combines 8-bit R/G/B components into the COM 0xBBGGRR integer
# format expected by ThemeColor.RGB.
function RGB([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

$tcs.Colors(1).RGB  = RGB 0x00 0x00 0x00   # dk1      -> 000000
$tcs.Colors(2).RGB  = RGB 0xFF 0xFF 0xFF   # lt1      -> FFFFFF
$tcs.Colors(3).RGB  = RGB 0x44 0x54 0x6A   # dk2      -> 44546A
$tcs.Colors(4).RGB  = RGB 0xE7 0xE6 0xE6   # lt2      -> E7E6E6
$tcs.Colors(5).RGB  = RGB 0x5B 0x9B 0xD5   # accent1  -> 5B9BD5
$tcs.Colors(6).RGB  = RGB 0xED 0x7D 0x31   # accent2  -> ED7D31
$tcs.Colors(7).RGB  = RGB 0xA5 0xA5 0xA5   # accent3  -> A5A5A5
$tcs.Colors(8).RGB  = RGB 0xFF 0xC0 0x00   # accent4  -> FFC000
$tcs.Colors(9).RGB  = RGB 0x44 0x72 0xC4   # accent5  -> 4472C4
$tcs.Colors(10).RGB = RGB 0x70 0xAD 0x47   # accent6  -> 70AD47
$tcs.Colors(11).RGB = RGB 0x05 0x63 0xC1   # hlink    -> 0563C1
$tcs.Colors(12).RGB = RGB 0x95 0x4F 0x72   # folHlink -> 954F72
